# Fix contact information missing from short resumes.
# Insert a new centered paragraph containing the contact info line
# directly after the "Dheeraj Chand" title paragraph (and before the
# "PROFESSIONAL SUMMARY" heading paragraph), matching the long-resume
# layout.

$d = $word.ActiveDocument

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

# Position right before the paragraph mark that ends the title
# paragraph ("Dheeraj Chand") so the new paragraph break inherits the
# title paragraph's formatting (centered), not the following
# paragraph's (Heading2).
$breakPoint = $d.Range($titleRange.End - 1, $titleRange.End - 1)
$breakPoint.InsertBefore("`r")

# The freshly created paragraph is now paragraph 2 - a blank, centered
# paragraph with no run formatting. Fill in the contact info text.
$contactPara = $d.Paragraphs.Item(2)
$contactPara.Range.InsertBefore("202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX")

Write-Output "Inserted contact info paragraph"
